$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E4").Value = '[5.00485667600443, 5.983461121429081, 6.314651948096729]'
$ws.Range("F4").Value = 380.103
$ws.Range("G4").Value = 380.103
$ws.Range("H4").Value = 19.4962
$ws.Range("I4").Value = 15.3183

# Row 7
$ws.Range("D7").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E7").Value = '[5.000999585973261, 6.279473648005702, 1.2259858884162598]'
$ws.Range("F7").Value = 410.43
$ws.Range("G7").Value = 410.43
$ws.Range("H7").Value = 20.2591
$ws.Range("I7").Value = 15.8307

# Row 10
$ws.Range("D10").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E10").Value = '[5.0182969558324535, 6.2831062525375065, 1.1693832191200728]'
$ws.Range("F10").Value = 15.821
$ws.Range("G10").Value = 411.899
$ws.Range("H10").Value = 20.2953
$ws.Range("I10").Value = 15.821
$ws.Range("J10").Value = 5142

# Row 13
$ws.Range("D13").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E13").Value = '[5.005316541017875, 6.238408919205355, 0.5592684043052151]'
$ws.Range("F13").Value = 15.9375
$ws.Range("G13").Value = 417.213
$ws.Range("H13").Value = 20.4258
$ws.Range("I13").Value = 15.9375
$ws.Range("J13").Value = 9563

# Row 16
$ws.Range("D16").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E16").Value = '[3.521400145956116, 0.5552473685876581, 0.05554223146572104]'
$ws.Range("F16").Value = 89.06059999999999
$ws.Range("G16").Value = 12324.8
$ws.Range("H16").Value = 111.017
$ws.Range("I16").Value = 89.06059999999999

# Row 19
$ws.Range("D19").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E19").Value = '[5.001810207951951, 6.217262432278078, 2.2954409742056643]'
$ws.Range("F19").Value = 20.0545
$ws.Range("G19").Value = 402.181
$ws.Range("H19").Value = 20.0545
$ws.Range("I19").Value = 15.682

# Row 22
$ws.Range("D22").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E22").Value = '[5.000017644910657, 6.2949540729247335, 0.5205321685482925]'
$ws.Range("F22").Value = 20.4086
$ws.Range("G22").Value = 416.513
$ws.Range("H22").Value = 20.4086
$ws.Range("I22").Value = 15.9358

# Row 25
$ws.Range("D25").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E25").Value = '[4.450693914280097, 0.652114566105611, 0.047605284169502285]'
$ws.Range("F25").Value = 80.2854
$ws.Range("G25").Value = 6445.75
$ws.Range("H25").Value = 80.2854
$ws.Range("I25").Value = 65.4198

# Row 28
$ws.Range("D28").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E28").Value = '[4.983957676271402, 5.739800903498828, 636.1702426335031]'
$ws.Range("F28").Value = 374.03
$ws.Range("G28").Value = 391954

# Row 31
$ws.Range("D31").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E31").Value = '[2.08918469041902, 3.2032997753214736, 288.41628469413746]'
$ws.Range("F31").Value = 29038.5
$ws.Range("G31").Value = 102072
$ws.Range("H31").Value = 319.488
$ws.Range("I31").Value = 276.283

# Row 34
$ws.Range("D34").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E34").Value = '[0.2653006148954005, 0.3946765474976315, 39.76616411296711]'
$ws.Range("F34").Value = 77748.39999999999
$ws.Range("G34").Value = 76911.39999999999
$ws.Range("H34").Value = 277.329
$ws.Range("I34").Value = 196.237

# Row 37
$ws.Range("D37").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E37").Value = '[0.03385930158884735, 0.0037733523816843536, 7.8661388107914885]'
$ws.Range("F37").Value = 213.972
$ws.Range("G37").Value = 85500.2
$ws.Range("H37").Value = 292.404

# Row 40
$ws.Range("D40").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E40").Value = '[0.0033329106634632055, 0.001012924300742405, 0.8102775227351402]'
$ws.Range("F40").Value = 216.419
$ws.Range("G40").Value = 86962.89999999999
$ws.Range("H40").Value = 294.895

# Row 43
$ws.Range("D43").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E43").Value = '[0.0002619441382201494, 0.001492323818604695, 0.081943086654686]'
$ws.Range("F43").Value = 216.68
$ws.Range("G43").Value = 87115.89999999999
$ws.Range("H43").Value = 295.154

# Row 46
$ws.Range("D46").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E46").Value = '[0.04617389420109781, 0.06901858239036192, 7.014823830265563]'
$ws.Range("F46").Value = 292.317
$ws.Range("G46").Value = 85058.60000000001
$ws.Range("H46").Value = 291.648

# Row 49
$ws.Range("D49").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E49").Value = '[0.00454611173967636, 0.008271186566630732, 0.7040332720277552]'
$ws.Range("F49").Value = 294.897
$ws.Range("G49").Value = 86919.8
$ws.Range("H49").Value = 294.822

# Row 52
$ws.Range("D52").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E52").Value = '[0.000383264869602737, 0.002217929497130577, 0.07133508842283416]'
$ws.Range("F52").Value = 295.155
$ws.Range("G52").Value = 87111.60000000001
$ws.Range("H52").Value = 295.147

# Row 55
$ws.Range("D55").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E55").Value = '[5.00783665597132, 5.73389279273554, 10.599069856939614]'

# Row 58
$ws.Range("D58").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E58").Value = '[4.330446172762307, 4.958621680177322, 9.163831844475599]'
$ws.Range("F58").Value = 1960.58
$ws.Range("G58").Value = 1960.58
$ws.Range("H58").Value = 44.2784
$ws.Range("I58").Value = 33.7473

# Row 61
$ws.Range("D61").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E61").Value = '[0.9076981297056335, 1.0412981435721893, 1.911798902988247]'
$ws.Range("F61").Value = 58530.9
$ws.Range("G61").Value = 58530.9
$ws.Range("H61").Value = 241.932
$ws.Range("I61").Value = 177.745

# Row 64
$ws.Range("D64").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E64").Value = '[0.11775335804600519, 0.1507718736850582, 2.864600516129187]'
$ws.Range("F64").Value = 210.632
$ws.Range("G64").Value = 82895.5
$ws.Range("H64").Value = 287.916
$ws.Range("I64").Value = 210.632

# Row 67
$ws.Range("D67").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E67").Value = '[0.011655677585664468, 0.016324290508674687, 0.30356799458475797]'
$ws.Range("F67").Value = 216.087
$ws.Range("G67").Value = 86704.3
$ws.Range("H67").Value = 294.456
$ws.Range("I67").Value = 216.087

# Row 70
$ws.Range("D70").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E70").Value = '[0.0010276986188542264, 0.0036712773259981344, 0.01979641454084068]'
$ws.Range("F70").Value = 216.647
$ws.Range("G70").Value = 87092
$ws.Range("H70").Value = 295.114
$ws.Range("I70").Value = 216.647

# Row 73
$ws.Range("D73").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E73").Value = '[0.16949820453171957, 0.1964308307502374, 0.34771936589672325]'
$ws.Range("F73").Value = 285.231
$ws.Range("G73").Value = 81356.60000000001
$ws.Range("H73").Value = 285.231
$ws.Range("I73").Value = 209.414

# Row 76
$ws.Range("D76").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E76").Value = '[0.01681309804808745, 0.021683243784050035, 0.024213981344683125]'
$ws.Range("F76").Value = 294.191
$ws.Range("G76").Value = 86548.3
$ws.Range("H76").Value = 294.191
$ws.Range("I76").Value = 215.979

# Row 79
$ws.Range("D79").Value = '[-0.008812772562144786, 0.019288968587611602, 0.01098815739164837]'
$ws.Range("E79").Value = '[0.0015434517499665515, 0.004207185340131793, -0.008138963296820134]'
$ws.Range("F79").Value = 295.087
$ws.Range("G79").Value = 87076.39999999999
$ws.Range("H79").Value = 295.087
$ws.Range("I79").Value = 216.637
